# FinalProjReflection.docx — "added premise and game description"
#
# The original document had a standalone paragraph containing only the
# text "Design". The edit inserts two new labeled sections ("Game
# Description:" and "Premise:"), each followed by an explanatory
# paragraph, directly above that "Design" paragraph — the "Design"
# paragraph itself is left untouched as the last paragraph of the group.

$d = $word.ActiveDocument

# Locate the paragraph that contains only "Design" (the anchor point).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n") -eq "Design") {
        $target = $p
        break
    }
}

$gameDescriptionLabel = "Game Description:"
$gameDescriptionBody = "The objective of the game is simple. Collect the three types of sushi: Tuna, Salmon, and Unagi. You must collect all Tuna first then get to the end and put all Tuna down on the plate. Then you can get all the Salmon and put it on the plate. Then you must get the final sushi the Unagi sushi and put it on the plate. Throughout your journey you will encounter puzzles that test your physical & mental strength along with your luck. You must successfully win these puzzles or games and receive the sushis. If you are successful and win the game, then you get to return back to your village with your new sushi skills and will be able to start your own sushi restaurant!  "

$premiseLabel = "Premise:"
$premiseBody = "Sushi is a rare and scarce commodity to the land. In fact it has been so sacred that the local government has banned sushi forever from all villages to keep to themselves. You are on an adventure to gather up sushi throughout many trials that test your physical strength, mental strength, and luck. Goal: Gather up the required sushi: tuna, salmon, and unagi, and you will be known as the sushi hero and will be able to go back to your land and start your very own Sushi restaurant at your village."

# Insert the four new paragraphs immediately before "Design" via
# InsertBefore so the new paragraphs inherit the same paragraph/run
# formatting (Helvetica, color 2D3B45, sz 24) as "Design". Each call
# inserts right in front of the (unmoved) "Design" paragraph, so build
# the block in reverse so the final reading order is correct:
#   Game Description: / <body> / Premise: / <body> / Design
$target.Range.InsertBefore($premiseBody + "`r")
$target.Range.InsertBefore($premiseLabel + "`r")
$target.Range.InsertBefore($gameDescriptionBody + "`r")
$target.Range.InsertBefore($gameDescriptionLabel + "`r")
